$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change
$ws.Range("B3").Value = 25034093703317.76
$ws.Range("C3").Value = 24103772653532.32
$ws.Range("D3").Value = 27063885084667.86

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 27942160201598.19
$ws.Range("C4").Value = 35012859283096.58
$ws.Range("D4").Value = 34391411772179.32

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 51825640668510.01
$ws.Range("C5").Value = 61019790202618.31
$ws.Range("D5").Value = 93612343722462.41
